# punto_3.xlsx: update the "coverage" column (E) with the final simulation
# results (re-run with a fixed seed, so the theta=5 rows now match theta=2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$coverage = @{
    2  = "0.955"
    3  = "0.9498"
    4  = "0.9498"
    5  = "0.9484"
    6  = "0.945"
    7  = "0.945"
    8  = "0.9468"
    9  = "0.9472"
    10 = "0.9472"
    11 = "0.9474"
    12 = "0.9464"
    13 = "0.9464"
    14 = "0.955"
    15 = "0.9498"
    16 = "0.9498"
    17 = "0.9484"
    18 = "0.945"
    19 = "0.945"
    20 = "0.9468"
    21 = "0.9472"
    22 = "0.9472"
    23 = "0.9474"
    24 = "0.9464"
    25 = "0.9464"
}

foreach ($row in $coverage.Keys) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.NumberFormat = "@"
    $cell.Value = $coverage[$row]
    $cell.Style = "Normal"
}

# Column E was manually narrowed a touch (no longer auto "best fit").
$ws.Columns.Item(5).ColumnWidth = 10.35

# Selection left on the whole table range instead of a stray cell.
$ws.Range("A1:E25").Select()
